$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddOpportunity")
$ws.Activate()
$ws.Range("AA2").Value = "10000.0"
$ws.Range("U19").Select()
